$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'68.120.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.86%  "

$ws.Range("D3").Formula = "'3.642.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.52%  "

$ws.Range("D4").Formula = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Formula = "'587.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("E6").Value = "  +5.15%  "

$ws.Range("D7").Formula = "'3.640.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.29%  "

$ws.Range("E8").Value = "  -5.97%  "

$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Formula = "'0.708"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.64%  "

$ws.Range("E11").Value = "  -9.58%  "

$ws.Range("D12").Formula = "'55.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.59%  "

$ws.Range("E13").Value = "  -10.52%  "

$ws.Range("D14").Formula = "'10.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.38%  "

$ws.Range("D15").Formula = "'4.224.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.77%  "

$ws.Range("D16").Formula = "'3.644.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.07%  "

$ws.Range("D17").Formula = "'19.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.23%  "

$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("E19").Value = "  -7.18%  "

$ws.Range("D20").Formula = "'12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.73%  "

$ws.Range("D21").Formula = "'67.957.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.49%  "

$ws.Range("D22").Formula = "'407.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.71%  "

$ws.Range("D23").Formula = "'4.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.54%  "

$ws.Range("D24").Formula = "'87.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.83%  "

$ws.Range("D25").Formula = "'2.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.05%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Formula = "'12.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.69%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Formula = "'3.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "

$ws.Range("E28").Value = "  -7.83%  "

$ws.Range("D29").Formula = "'6.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").Formula = "'9.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.45%  "

$ws.Range("D31").Formula = "'32.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.37%  "

$ws.Range("D32").Formula = "'7.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.77%  "

$ws.Range("D33").Formula = "'12.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.96%  "

$ws.Range("E34").Value = "  -7.34%  "

$ws.Range("D35").Formula = "'64.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.56%  "

$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Formula = "'599.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.43%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Formula = "'42.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.94%  "

$ws.Range("D38").Formula = "'0.0₃0878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.22%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Formula = "'0.394"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.38%  "

$ws.Range("D41").Formula = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").Formula = "'0.136"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.94%  "

$ws.Range("D43").Formula = "'3.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.22%  "

$ws.Range("D44").Formula = "'2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.23%  "

$ws.Range("D45").Formula = "'0.0435"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.28%  "

$ws.Range("D46").Formula = "'2.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -13.09%  "

$ws.Range("D47").Formula = "'0.133"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.63%  "

$ws.Range("D48").Formula = "'2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("D49").Formula = "'8.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.40%  "

$ws.Range("D50").Formula = "'3.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.99%  "

$ws.Range("D51").Formula = "'2.685.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.02%  "
